# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the Slide Master / all slides)
#
# The authored edit swaps the two themes' color palettes so the slides
# (theme2.xml) now render with the default "Office Theme" color scheme
# instead of "Integral". Re-color every theme-color slot on the slide
# master via the ThemeColorScheme that PowerPoint exposes on a Slide
# (it edits the shared master theme, not a per-slide override) using the
# stock Office theme RGB values, in msoThemeColorIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hyperlink, 12 followedHyperlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(3).RGB  = 6968388    # dk2       44546A
$tcs.Item(4).RGB  = 15132391   # lt2       E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1   5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2   ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3   A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4   FFC000
$tcs.Item(9).RGB  = 12874308   # accent5   4472C4
$tcs.Item(10).RGB = 4697456    # accent6   70AD47
$tcs.Item(11).RGB = 12673797   # hlink     0563C1
$tcs.Item(12).RGB = 7491477    # folHlink  954F72
